$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header cells (these are the Tabla1 structured-table headers,
# so renaming the header cells renames the table columns too):
#   C1: "Descripción (nombre producto/servicio)" -> "Nombre Producto/Servicio"
#   D1: "Nombre Etiqueta (quitar para bsale)"     -> "Nombre Etiqueta"
# C1 is renamed before D1 (matches shared-string slot reuse order), and both
# happen before the formula loop below, which needs Tabla1[...[Nombre Etiqueta]]
# to resolve against the new column name.
$ws.Range("C1").Value = "Nombre Producto/Servicio"
$ws.Range("D1").Value = "Nombre Etiqueta"

# Add the label-concatenation formula down column C for every data row (3-96).
# Row 2 already carries a literal value for C2 and is left untouched.
# Cells are written one at a time (not as a single multi-cell range) so each
# one gets its own full, non-shared formula definition.
$formula = '=CONCATENATE(Tabla1[[#This Row],[Nombre Etiqueta]]," ",Tabla1[[#This Row],[Tamanio]]," ",Tabla1[[#This Row],[Posicion]])'
for ($r = 3; $r -le 96; $r++) {
    $ws.Cells.Item($r, 3).Formula = $formula
}

# Narrow column C slightly now that its header text is shorter.
$ws.Columns.Item(3).ColumnWidth = 34

# Move the visible selection to D8 (matches the saved cursor position).
$ws.Range("D8").Select()
